$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new weekly price row above the current row 166 ("Feria Lagunitas
# de Puerto Montt" - Cebollín table). This shifts existing rows 166-208 down
# to 167-209, extending the sheet's used range to A1:R209.
$ws.Rows("166:166").Insert()

# Populate the newly inserted row 166 with the new observation.
$ws.Cells.Item(166, 1).Value = 4
$ws.Cells.Item(166, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(166, 3).Value = "Los Lagos"
$ws.Cells.Item(166, 4).Value = 44551
$ws.Cells.Item(166, 5).Value = 10
$ws.Cells.Item(166, 6).Value = 100112037
$ws.Cells.Item(166, 7).Value = "Cebollín"
$ws.Cells.Item(166, 8).Value = "Sin especificar"
$ws.Cells.Item(166, 9).Value = "Primera"
$ws.Cells.Item(166, 10).Value = 180
$ws.Cells.Item(166, 11).Value = 7000
$ws.Cells.Item(166, 12).Value = 7000
$ws.Cells.Item(166, 13).Value = 7000
$ws.Cells.Item(166, 14).Value = "`$/paquete 36 unidades"
$ws.Cells.Item(166, 15).Value = "Región Metropolitana"
$ws.Cells.Item(166, 16).Value = 194
$ws.Cells.Item(166, 17).Value = 36
$ws.Cells.Item(166, 18).Value = "Hortaliza"
